# MOS-23045: Update Master Data as per 22 April Changes
# Adds 10 new "Postal Code" location rows (eng/fra/ara) under parent BNMR
# to the master-location sheet (rows 110-119).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(110, 10110, 10110, 5, "Postal Code", "BNMR", "eng"),
    @(111, 10111, 10111, 5, "Postal Code", "BNMR", "eng"),
    @(112, 10113, 10113, 5, "Postal Code", "BNMR", "eng"),
    @(113, 10114, 10114, 5, "Postal Code", "BNMR", "eng"),
    @(114, 10111, 10111, 5, "code postal", "BNMR", "fra"),
    @(115, 10110, 10110, 5, "code postal", "BNMR", "fra"),
    @(116, 10113, 10113, 5, "code postal", "BNMR", "fra"),
    @(117, 10114, 10114, 5, "code postal", "BNMR", "fra"),
    @(118, 10111, 10111, 5, "الرمز البريدي", "BNMR", "ara"),
    @(119, 10110, 10110, 5, "الرمز البريدي", "BNMR", "ara")
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $ws.Cells.Item($rowIndex, 1).Value = $r[1]
    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]
    $ws.Cells.Item($rowIndex, 5).Value = $r[5]
    $ws.Cells.Item($rowIndex, 6).Value = $r[6]
    $ws.Cells.Item($rowIndex, 7).Value = $true
    $ws.Cells.Item($rowIndex, 8).Value = "superadmin"
    $ws.Cells.Item($rowIndex, 9).Value = "now()"
}
